$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: CageID for row 2 was wrong (12A -> 50A)
$ws.Range("D2").Value = "50A"

# Row 12
$ws.Range("A12").Value = 54332
$ws.Range("B12").Value = "American Gouldian"
$ws.Range("C12").Value = "North America"
$ws.Range("D12").Value = "22A"
$ws.Range("E12").Value = "Male"
$ws.Range("F12").Value = 223
$ws.Range("G12").Value = 111
$ws.Range("H12").Value = "15/05/2023"
$ws.Range("I12").Value = "Red"
$ws.Range("J12").Value = "Purple"
$ws.Range("K12").Value = "Pastel"

# Row 13
$ws.Range("A13").Value = 543322
$ws.Range("B13").Value = "American Gouldian"
$ws.Range("C13").Value = "North America"
$ws.Range("D13").Value = "22A"
$ws.Range("E13").Value = "Female"
$ws.Range("F13").Value = 223
$ws.Range("G13").Value = 111
$ws.Range("H13").Value = "15/05/2023"
$ws.Range("I13").Value = "Black"
$ws.Range("J13").Value = "Purple"
$ws.Range("K13").Value = "Green"

# Row 14
$ws.Range("A14").Value = 3425
$ws.Range("B14").Value = "American Gouldian"
$ws.Range("C14").Value = "Central America"
$ws.Range("D14").Value = "50A"
$ws.Range("E14").Value = "Male"
$ws.Range("F14").Value = 101
$ws.Range("G14").Value = 1231
$ws.Range("H14").Value = "18/05/2023"
$ws.Range("I14").Value = "Red"
$ws.Range("J14").Value = "Purple"
$ws.Range("K14").Value = "Green Pastel"

# Row 15
$ws.Range("A15").Value = 222
$ws.Range("B15").Value = "American Gouldian"
$ws.Range("C15").Value = "Central America"
$ws.Range("D15").Value = 444
$ws.Range("E15").Value = "Female"
$ws.Range("F15").Value = 101
$ws.Range("G15").Value = 1231
$ws.Range("H15").Value = 45144
$ws.Range("H2").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("I15").Value = "Red"
$ws.Range("J15").Value = "Purple"
$ws.Range("K15").Value = "Green Pastel"
